$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Checklist" to "Session"
$ws.Name = "Session"

# Row 2 is updated to the values previously on row 3 (student id 191216,
# log time 14:14:24), with the log Type switched from "Selection" to "Scan".
# Force text format on the two columns Excel would otherwise auto-coerce
# (A: numeric-looking ID, C: date-looking string) so they stay text, matching
# the sheet's existing "numberStoredAsText" handling.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"

$ws.Range("A2").Value = "191216"
$ws.Range("B2").Value = "Pediatrics"
$ws.Range("C2").Value = "06/09/2025"
$ws.Range("D2").Value = "14:14:24"
$ws.Range("E2").Value = "Scan"
$ws.Range("F2").Value = "user@user.com"

# The old row 3 (now duplicated into row 2) is removed entirely, shifting
# rows up so the sheet ends at row 2.
$ws.Rows.Item(3).Delete()
